$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures from the latest GitHub Actions run.
# Each D/E cell holds numeric-looking text (e.g. "300.85", "-0.81%"); the source sheet
# stores these as literal text, so force Text number format before assigning the value
# to stop Excel auto-coercing them into floating point numbers / percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.81%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.01%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.114"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.84%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07360"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.90%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.320"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "55.86%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.949"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.17%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.73%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9184"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.29%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1704"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.17%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07614"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08101"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.58%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02985"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09926"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.23%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.26%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006183"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.80%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.13%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3292"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.31%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1319"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.27%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.655"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.37%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.71%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.24%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.85%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004479"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.72%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.06%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-3.17%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01729"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.60%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04509"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007245"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.10%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1343"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.37%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01072"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.91%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006273"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.82%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-33.27%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "13.69%"
